# database_glossary.xlsx — "Add error-catching where data is missing"
#
# The sheet gains a leading spacer row + spacer column (so the table starts
# at B2 instead of A1), a couple of helper cells in column G that
# cross-check the incentive cap (W_incentive_max_USD / incentive_per_W),
# and a red-on-pink "needs attention" highlight on the rows whose source
# data is most likely to go missing/be stale.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Incentive testing")

# --- Shift the whole table one column right and one row down ------------
# (Excel re-addresses every formula/range automatically, same as a user
# choosing Home > Insert > Insert Sheet Columns / Insert Sheet Rows.)
$ws1.Columns("A:A").Insert()
$ws1.Rows("1:1").Insert()

# --- New helper cells in column G ----------------------------------------
# Row 9  = electricity_price_old, row 10 = electricity_price (post-shift)
$ws1.Range("G9").Formula = "=28296/20"

$ws1.Range("C10").Copy()
$ws1.Range("G10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws1.Range("G10").Value = 5000

$ws1.Range("G11").Formula = "=G9/G10"

# --- Highlight the rows most likely to be missing/incomplete data --------
# Red font on the existing pink fill (fillId 3) already used on Sheet2's
# header cells — pull that exact fill across via PasteSpecial so no new
# fill entry is created, then tint the font red.
$flagged = @("B10", "B30", "B31", "B32", "B33", "B36")
$ws2.Range("E1").Copy()
foreach ($addr in $flagged) {
    $ws1.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false
foreach ($addr in $flagged) {
    $ws1.Range($addr).Font.Color = 255
}

# --- View tidy-up: zoom out a touch and rest the selection on B2 ---------
$ws1.Activate()
$excel.ActiveWindow.Zoom = 144
$ws1.Range("B2").Select()
